# Natmi re-run (Hgf -> St14) following Dr Hou's advice: the sending/target
# cluster set now includes "ECs", giving a 4x4 grid (16 data rows, up from
# 3x4/12) with refreshed expression/specificity statistics for every pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "St14"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 8.265822
$ws.Range("H2").Value = 24.797466
$ws.Range("I2").Value = 0.2082338764513023
$ws.Range("J2").Value = 0.2082338764513023
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7597586666666668
$ws.Range("N2").Value = 2.279276
$ws.Range("O2").Value = 0.3177111952105157
$ws.Range("P2").Value = 0.3177111952105158
$ws.Range("Q2").Value = 6.280029901624001
$ws.Range("R2").Value = 56.52026911461601
$ws.Range("S2").Value = 0.0661582337706621
$ws.Range("T2").Value = 0.06615823377066211

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "St14"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 8.265822
$ws.Range("H3").Value = 24.797466
$ws.Range("I3").Value = 0.2082338764513023
$ws.Range("J3").Value = 0.2082338764513023
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6247346666666667
$ws.Range("N3").Value = 1.874204
$ws.Range("O3").Value = 0.2612476913319534
$ws.Range("P3").Value = 0.2612476913319534
$ws.Range("Q3").Value = 5.163945551896
$ws.Range("R3").Value = 46.475509967064
$ws.Range("S3").Value = 0.05440061948000592
$ws.Range("T3").Value = 0.05440061948000592

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "St14"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 8.265822
$ws.Range("H4").Value = 24.797466
$ws.Range("I4").Value = 0.2082338764513023
$ws.Range("J4").Value = 0.2082338764513023
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3417453333333333
$ws.Range("N4").Value = 1.025236
$ws.Range("O4").Value = 0.1429089565865864
$ws.Range("P4").Value = 0.1429089565865864
$ws.Range("Q4").Value = 2.824806094664
$ws.Range("R4").Value = 25.423254851976
$ws.Range("S4").Value = 0.02975848600963574
$ws.Range("T4").Value = 0.02975848600963575

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "St14"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 8.265822
$ws.Range("H5").Value = 24.797466
$ws.Range("I5").Value = 0.2082338764513023
$ws.Range("J5").Value = 0.2082338764513023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6651113333333333
$ws.Range("N5").Value = 1.995334
$ws.Range("O5").Value = 0.2781321568709446
$ws.Range("P5").Value = 0.2781321568709446
$ws.Range("Q5").Value = 5.497691891515999
$ws.Range("R5").Value = 49.47922702364399
$ws.Range("S5").Value = 0.05791653719099848
$ws.Range("T5").Value = 0.05791653719099848

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "St14"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.28595333333333
$ws.Range("H6").Value = 33.85786
$ws.Range("I6").Value = 0.2843174958338682
$ws.Range("J6").Value = 0.2843174958338682
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7597586666666668
$ws.Range("N6").Value = 2.279276
$ws.Range("O6").Value = 0.3177111952105157
$ws.Range("P6").Value = 0.3177111952105158
$ws.Range("Q6").Value = 8.574600856595557
$ws.Range("R6").Value = 77.17140770936001
$ws.Range("S6").Value = 0.0903308514206391
$ws.Range("T6").Value = 0.09033085142063911

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "St14"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.28595333333333
$ws.Range("H7").Value = 33.85786
$ws.Range("I7").Value = 0.2843174958338682
$ws.Range("J7").Value = 0.2843174958338682
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6247346666666667
$ws.Range("N7").Value = 1.874204
$ws.Range("O7").Value = 0.2612476913319534
$ws.Range("P7").Value = 0.2612476913319534
$ws.Range("Q7").Value = 7.050726293715556
$ws.Range("R7").Value = 63.45653664344
$ws.Range("S7").Value = 0.07427728939188034
$ws.Range("T7").Value = 0.07427728939188034

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "St14"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.28595333333333
$ws.Range("H8").Value = 33.85786
$ws.Range("I8").Value = 0.2843174958338682
$ws.Range("J8").Value = 0.2843174958338682
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3417453333333333
$ws.Range("N8").Value = 1.025236
$ws.Range("O8").Value = 0.1429089565865864
$ws.Range("P8").Value = 0.1429089565865864
$ws.Range("Q8").Value = 3.856921883884445
$ws.Range("R8").Value = 34.71229695496
$ws.Range("S8").Value = 0.04063151666892922
$ws.Range("T8").Value = 0.04063151666892923

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "St14"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.28595333333333
$ws.Range("H9").Value = 33.85786
$ws.Range("I9").Value = 0.2843174958338682
$ws.Range("J9").Value = 0.2843174958338682
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6651113333333333
$ws.Range("N9").Value = 1.995334
$ws.Range("O9").Value = 0.2781321568709446
$ws.Range("P9").Value = 0.2781321568709446
$ws.Range("Q9").Value = 7.50641546947111
$ws.Range("R9").Value = 67.55773922524
$ws.Range("S9").Value = 0.07907783835241956
$ws.Range("T9").Value = 0.07907783835241956

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "St14"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.51551966666667
$ws.Range("H10").Value = 58.546559
$ws.Range("I10").Value = 0.4916380138783083
$ws.Range("J10").Value = 0.4916380138783083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7597586666666668
$ws.Range("N10").Value = 2.279276
$ws.Range("O10").Value = 0.3177111952105157
$ws.Range("P10").Value = 0.3177111952105158
$ws.Range("Q10").Value = 14.82708520125378
$ws.Range("R10").Value = 133.443766811284
$ws.Range("S10").Value = 0.1561989010002015
$ws.Range("T10").Value = 0.1561989010002015

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Hgf"
$ws.Range("C11").Value = "St14"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 19.51551966666667
$ws.Range("H11").Value = 58.546559
$ws.Range("I11").Value = 0.4916380138783083
$ws.Range("J11").Value = 0.4916380138783083
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6247346666666667
$ws.Range("N11").Value = 1.874204
$ws.Range("O11").Value = 0.2612476913319534
$ws.Range("P11").Value = 0.2612476913319534
$ws.Range("Q11").Value = 12.19202167378178
$ws.Range("R11").Value = 109.728195064036
$ws.Range("S11").Value = 0.1284392960967349
$ws.Range("T11").Value = 0.1284392960967349

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Hgf"
$ws.Range("C12").Value = "St14"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 19.51551966666667
$ws.Range("H12").Value = 58.546559
$ws.Range("I12").Value = 0.4916380138783083
$ws.Range("J12").Value = 0.4916380138783083
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.3417453333333333
$ws.Range("N12").Value = 1.025236
$ws.Range("O12").Value = 0.1429089565865864
$ws.Range("P12").Value = 0.1429089565865864
$ws.Range("Q12").Value = 6.669337773658222
$ws.Range("R12").Value = 60.02403996292401
$ws.Range("S12").Value = 0.0702594755816507
$ws.Range("T12").Value = 0.07025947558165072

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Hgf"
$ws.Range("C13").Value = "St14"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 19.51551966666667
$ws.Range("H13").Value = 58.546559
$ws.Range("I13").Value = 0.4916380138783083
$ws.Range("J13").Value = 0.4916380138783083
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6651113333333333
$ws.Range("N13").Value = 1.995334
$ws.Range("O13").Value = 0.2781321568709446
$ws.Range("P13").Value = 0.2781321568709446
$ws.Range("Q13").Value = 12.97999330618955
$ws.Range("R13").Value = 116.819939755706
$ws.Range("S13").Value = 0.1367403411997213
$ws.Range("T13").Value = 0.1367403411997213

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Hgf"
$ws.Range("C14").Value = "St14"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6276006666666666
$ws.Range("H14").Value = 1.882802
$ws.Range("I14").Value = 0.01581061383652123
$ws.Range("J14").Value = 0.01581061383652123
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.7597586666666668
$ws.Range("N14").Value = 2.279276
$ws.Range("O14").Value = 0.3177111952105157
$ws.Range("P14").Value = 0.3177111952105158
$ws.Range("Q14").Value = 0.4768250457057778
$ws.Range("R14").Value = 4.291425411352
$ws.Range("S14").Value = 0.005023209019013077
$ws.Range("T14").Value = 0.005023209019013077

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Hgf"
$ws.Range("C15").Value = "St14"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6276006666666666
$ws.Range("H15").Value = 1.882802
$ws.Range("I15").Value = 0.01581061383652123
$ws.Range("J15").Value = 0.01581061383652123
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.6247346666666667
$ws.Range("N15").Value = 1.874204
$ws.Range("O15").Value = 0.2612476913319534
$ws.Range("P15").Value = 0.2612476913319534
$ws.Range("Q15").Value = 0.3920838932897777
$ws.Range("R15").Value = 3.528755039608
$ws.Range("S15").Value = 0.004130486363332209
$ws.Range("T15").Value = 0.004130486363332209

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Hgf"
$ws.Range("C16").Value = "St14"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6276006666666666
$ws.Range("H16").Value = 1.882802
$ws.Range("I16").Value = 0.01581061383652123
$ws.Range("J16").Value = 0.01581061383652123
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3417453333333333
$ws.Range("N16").Value = 1.025236
$ws.Range("O16").Value = 0.1429089565865864
$ws.Range("P16").Value = 0.1429089565865864
$ws.Range("Q16").Value = 0.2144795990302222
$ws.Range("R16").Value = 1.930316391272
$ws.Range("S16").Value = 0.002259478326370694
$ws.Range("T16").Value = 0.002259478326370694

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Hgf"
$ws.Range("C17").Value = "St14"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6276006666666666
$ws.Range("H17").Value = 1.882802
$ws.Range("I17").Value = 0.01581061383652123
$ws.Range("J17").Value = 0.01581061383652123
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6651113333333333
$ws.Range("N17").Value = 1.995334
$ws.Range("O17").Value = 0.2781321568709446
$ws.Range("P17").Value = 0.2781321568709446
$ws.Range("Q17").Value = 0.4174243162075555
$ws.Range("R17").Value = 3.756818845867999
$ws.Range("S17").Value = 0.004397440127805249
$ws.Range("T17").Value = 0.004397440127805249
